$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Main table A1:D6 gets a thin border around every cell ---
$ws.Range("A1:D6").Borders.LineStyle = 1

# --- New "Dado1"/"Dado2" header cells (K1:L1), bordered + centered ---
$ws.Range("K1").Value = "Dado1"
$ws.Range("L1").Value = "Dado2"
$ws.Range("K1:L1").Borders.LineStyle = 1
$ws.Range("K1:L1").HorizontalAlignment = -4108

# --- Helper label in J2 (plain, unstyled) ---
$ws.Range("J2").Value = " "

# --- Underlined marker cell below the table ---
$ws.Range("J8").Font.Underline = 2

# --- Merge + format the two "Leitura coluna" result columns ---
$ws.Range("K2:K6").Merge()
$ws.Range("K2").Value = 'Leitura coluna "k"'
$ws.Range("K2:K6").Borders.LineStyle = 1
$ws.Range("K2:K6").HorizontalAlignment = -4108
$ws.Range("K2:K6").VerticalAlignment = -4108

$ws.Range("L2:L6").Merge()
$ws.Range("L2").Value = 'Leitura coluna "l"'
$ws.Range("L2:L6").Borders.LineStyle = 1
$ws.Range("L2:L6").HorizontalAlignment = -4108
$ws.Range("L2:L6").VerticalAlignment = -4108

# --- Column widths for the new columns (auto-sized to fit content) ---
$ws.Columns("K:L").AutoFit()

# --- Selection / active cell matches the author's final view ---
$ws.Range("K1").Select()
